$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: new match data (Schalke vs Darmstadt, 24/10/2025)
$ws.Range("A11").Value = "24/10/2025"
$ws.Range("B11").Value = "Schalke"
$ws.Range("C11").Value = 1
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = "Darmstadt"
$ws.Range("F11").Value = "L"
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 1
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 0.45
$ws.Range("L11").Value = 1.84
$ws.Range("M11").Value = 7
$ws.Range("N11").Value = 10
$ws.Range("O11").Value = 1
$ws.Range("P11").Value = 2
